# Update column F (dSF) values on the worksheet to re-pulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = -3
$ws.Range("F4").Value  = -4
$ws.Range("F5").Value  = 4
$ws.Range("F6").Value  = -4
$ws.Range("F8").Value  = 4
$ws.Range("F9").Value  = -3
$ws.Range("F11").Value = -6
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 3
$ws.Range("F19").Value = -5
